$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename header C1 from "RunType" to "ConditionType"
$ws.Range("C1").Value = "ConditionType"

# Update Question (B) and RunType/ConditionType (C) columns for rows 2-20
$values = @(
    @(38, 32, 4),
    @(14, 33, 4),
    @(4, 34, 4),
    @(15, 24, 4),
    @(3, 18, 4),
    @(5, 37, 4),
    @(8, 25, 4),
    @(34, 10, 4),
    @(10, 30, 4),
    @(21, 13, 4),
    @(9, 9, 4),
    @(22, 4, 4),
    @(36, 6, 4),
    @(2, 26, 4),
    @(37, 2, 4),
    @(24, 15, 4),
    @(23, 38, 4),
    @(33, 36, 4),
    @(18, 21, 4)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $newB = $values[$i][1]
    $newC = $values[$i][2]
    $ws.Cells.Item($row, 2).Value = $newB
    $ws.Cells.Item($row, 3).Value = $newC
}
